$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1089366021538427
$ws.Range("B3").Value = 0.08500062045203191
$ws.Range("H3").Value = 0.1939372226058746
$ws.Range("B4").Value = 0.09917304706097013
$ws.Range("H4").Value = 0.2081096492148129
$ws.Range("B5").Value = 0.123807232963945
$ws.Range("C5").Value = 0.01369461941555491
$ws.Range("D5").Value = 16.00202245006286
$ws.Range("E5").Value = 0.06129388298243208
$ws.Range("F5").Value = 0.09691080096350797
$ws.Range("G5").Value = 0.1507036649643816
$ws.Range("H5").Value = 0.2327438351177878
$ws.Range("B6").Value = 0.1110720067567618
$ws.Range("H6").Value = 0.2200086089106045
$ws.Range("B7").Value = 0.03698560271094231
$ws.Range("H7").Value = 0.145922204864785
$ws.Range("B8").Value = 0.03114616664858009
$ws.Range("H8").Value = 0.1400827688024228
$ws.Range("B9").Value = 0.01659237862443653
$ws.Range("C9").Value = 0.002386670275222697
$ws.Range("D9").Value = 1.346453958426267
$ws.Range("E9").Value = 0.01391397053886815
$ws.Range("F9").Value = 0.01187945893565428
$ws.Range("G9").Value = 0.02130529831321876
$ws.Range("H9").Value = 0.1255289807782793
$ws.Range("B10").Value = 0.01571141166047701
$ws.Range("C10").Value = 0.001858316576873916
$ws.Range("D10").Value = 0.9755545145950986
$ws.Range("E10").Value = 0.003225062826648291
$ws.Range("F10").Value = 0.01205509558667424
$ws.Range("G10").Value = 0.01936772773428012
$ws.Range("H10").Value = 0.1246480138143198
$ws.Range("B11").Value = 0.02662391629174165
$ws.Range("H11").Value = 0.1355605184455844
$ws.Range("B12").Value = 0.038878014885981
$ws.Range("H12").Value = 0.1478146170398237
$ws.Range("B13").Value = 0.05065715289869716
$ws.Range("H13").Value = 0.1595937550525399
$ws.Range("B14").Value = 0.05419193279007294
$ws.Range("H14").Value = 0.1631285349439157
$ws.Range("B15").Value = 0.05995336894647041
$ws.Range("C15").Value = 0.008943036960297145
$ws.Range("D15").Value = 11.68031312222276
$ws.Range("E15").Value = 0.04716771729463884
$ws.Range("F15").Value = 0.04239731264739489
$ws.Range("G15").Value = 0.07750942524554597
$ws.Range("H15").Value = 0.1688899711003131
$ws.Range("B16").Value = 0.0642453643400073
$ws.Range("H16").Value = 0.17318196649385
$ws.Range("B17").Value = 0.06825029440056077
$ws.Range("C17").Value = 0.008841851205278278
$ws.Range("D17").Value = 12.66402477271134
$ws.Range("E17").Value = 0.05173022530280881
$ws.Range("F17").Value = 0.05089728137956876
$ws.Range("G17").Value = 0.08560330742155288
$ws.Range("H17").Value = 0.1771868965544035
$ws.Range("B18").Value = -0.1089366021538427
$ws.Range("C18").Value = 0.01281683042101939
$ws.Range("D18").Value = -16.44676351054337
$ws.Range("E18").Value = 0.04008260604847166
$ws.Range("F18").Value = -0.1341034378989942
$ws.Range("G18").Value = -0.08376976640869119
$ws.Range("B19").Value = 0.06830557377428821
$ws.Range("H19").Value = 0.177242175928131
$ws.Range("B20").Value = 0.07226336649908023
$ws.Range("C20").Value = 0.00890602765187856
$ws.Range("D20").Value = 5668150894494.086
$ws.Range("E20").Value = 0.04432099779570591
$ws.Range("F20").Value = 0.05478289827742396
$ws.Range("G20").Value = 0.08974383472073651
$ws.Range("H20").Value = 0.181199968652923
$ws.Range("B21").Value = 0.07241528648047106
$ws.Range("C21").Value = 0.008991318092873422
$ws.Range("D21").Value = -428133630578.8991
$ws.Range("E21").Value = 0.0506141623508175
$ws.Range("F21").Value = 0.05476696315514627
$ws.Range("G21").Value = 0.09006360980579604
$ws.Range("H21").Value = 0.1813518886343138
$ws.Range("B22").Value = 0.0729059892368354
$ws.Range("C22").Value = 0.008834673550572545
$ws.Range("D22").Value = 104298678.5056185
$ws.Range("E22").Value = 0.05542860158423735
$ws.Range("F22").Value = 0.0555645292479386
$ws.Range("G22").Value = 0.09024744922573248
$ws.Range("H22").Value = 0.1818425913906782
$ws.Range("B23").Value = 0.0738015138916179
$ws.Range("C23").Value = 0.00887600175227507
$ws.Range("D23").Value = 609921337538.8508
$ws.Range("E23").Value = 0.05043224462281864
$ws.Range("F23").Value = 0.05637587448914926
$ws.Range("G23").Value = 0.09122715329408662
$ws.Range("H23").Value = 0.1827381160454606
$ws.Range("B24").Value = 0.07547225321043659
$ws.Range("C24").Value = 0.008532349552489216
$ws.Range("D24").Value = 531793911756.0118
$ws.Range("E24").Value = 0.05373465134791208
$ws.Range("F24").Value = 0.05872494904565367
$ws.Range("G24").Value = 0.0922195573752197
$ws.Range("H24").Value = 0.1844088553642793
$ws.Range("B25").Value = 0.07771443315715504
$ws.Range("C25").Value = 0.009112689933835834
$ws.Range("D25").Value = 6634860042354.429
$ws.Range("E25").Value = 0.06172046840013464
$ws.Range("F25").Value = 0.05982240874266918
$ws.Range("G25").Value = 0.09560645757164075
$ws.Range("H25").Value = 0.1866510353109978
$ws.Range("B26").Value = 0.08000025974785395
$ws.Range("C26").Value = 0.009167721000248905
$ws.Range("D26").Value = 12.72034225711517
$ws.Range("E26").Value = 0.05365001257679801
$ws.Range("F26").Value = 0.06200389481192491
$ws.Range("G26").Value = 0.09799662468378292
$ws.Range("H26").Value = 0.1889368619016967
$ws.Range("B27").Value = 0.07900893735709619
$ws.Range("C27").Value = 0.009556243819865662
$ws.Range("D27").Value = 12.45274135783254
$ws.Range("E27").Value = 0.05898925413476505
$ws.Range("F27").Value = 0.06024596718288395
$ws.Range("G27").Value = 0.09777190753130867
$ws.Range("H27").Value = 0.1879455395109389
$ws.Range("B28").Value = 0.08342092917406124
$ws.Range("C28").Value = 0.009432738458176711
$ws.Range("D28").Value = 12.37454337021428
$ws.Range("E28").Value = 0.06874796859422674
$ws.Range("F28").Value = 0.06491092979499977
$ws.Range("G28").Value = 0.1019309285531226
$ws.Range("H28").Value = 0.192357531327904
$ws.Range("B29").Value = 0.008933096119192879
$ws.Range("C29").Value = 0.003874225444742033
$ws.Range("D29").Value = 1.163037890108816
$ws.Range("E29").Value = 0.03203492588249782
$ws.Range("F29").Value = 0.001305982207810635
$ws.Range("G29").Value = 0.01656021003057486
$ws.Range("H29").Value = 0.1178696982730356

Write-Output "Applied 134 cell updates"
